$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 fresh rows before the "ping" Spreadsheet block (old row 27) to make
# room for two new Spreadsheet method blocks: parse4QueryParam and parse4PathParam.
$ws.Rows("27:32").Insert()

# The insert copies formatting (quotePrefix style) from the row above into the
# blank rows - clear that so the new rows start out completely empty, matching
# the sparse row layout used throughout this sheet.
$ws.Range("C27:E32").ClearFormats()
$ws.Range("C27:E32").ClearContents()

# New Spreadsheet method: parse4QueryParam - same body shape as parse4.
$ws.Range("C28").Value = "Spreadsheet Integer parse4QueryParam(String str)"
$ws.Range("C29").Value = "Step"
$ws.Range("D29").Value = "Calc"
$ws.Range("C30").Value = "RETURN"
$ws.Range("D30").Value = "'= parse(str) + 100"

# New Spreadsheet method: parse4PathParam - same body shape as parse4.
$ws.Range("C32").Value = "Spreadsheet Integer parse4PathParam(String str)"
$ws.Range("C33").Value = "Step"
$ws.Range("D33").Value = "Calc"
$ws.Range("C34").Value = "RETURN"
$ws.Range("D34").Value = "'= parse(str) + 100"

# Reflect the author's final cursor position in the saved view state.
$ws.Range("C28").Select()
